$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -22.1153
$ws.Range("B4").Value = 5.841400000000002
$ws.Range("A7").Value = -20.06129999999997
$ws.Range("C10").Value = -13.23749999999999
$ws.Range("B12").Value = 4.842299999999998
$ws.Range("C13").Value = -13.6721
$ws.Range("A16").Value = -22.04690000000002
$ws.Range("B18").Value = 6.087199999999995
$ws.Range("B19").Value = 8.820600000000002
$ws.Range("B20").Value = 9.351299999999988
$ws.Range("D25").Value = -7.430999999999999
$ws.Range("A28").Value = -21.8641
$ws.Range("A29").Value = -21.15149999999997
$ws.Range("C30").Value = -12.24219999999999
$ws.Range("B31").Value = 4.8505
$ws.Range("A32").Value = -21.3186
$ws.Range("D34").Value = -7.853900000000005
$ws.Range("D39").Value = -8.323299999999994
$ws.Range("A40").Value = -21.86719999999998
$ws.Range("B40").Value = 6.077400000000003
$ws.Range("C40").Value = -12.1875
$ws.Range("B42").Value = 8.935999999999996
$ws.Range("C44").Value = -12.885
$ws.Range("B47").Value = 5.315400000000002
$ws.Range("B48").Value = 5.527100000000003
$ws.Range("A52").Value = -22.1361
$ws.Range("A57").Value = -22.28370000000001
$ws.Range("D61").Value = -8.150999999999996
$ws.Range("B63").Value = 4.837399999999999
$ws.Range("B64").Value = 5.513199999999999
$ws.Range("D64").Value = -7.18939999999999
$ws.Range("A66").Value = -21.45729999999999
$ws.Range("B76").Value = 5.690999999999997
$ws.Range("D78").Value = -7.948200000000003
$ws.Range("B81").Value = 4.855700000000004
$ws.Range("D83").Value = -8.600999999999996
$ws.Range("B89").Value = 4.934699999999996
$ws.Range("C89").Value = -13.3201
$ws.Range("C91").Value = -12.4735
$ws.Range("D92").Value = -6.344900000000001
$ws.Range("B94").Value = 4.677599999999995
$ws.Range("D98").Value = -7.602399999999998
$ws.Range("A100").Value = -22.07040000000002
